$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 7 : 2005 Dom Perignon "Brut Vintage" (Champagne) ----
$ws.Range("A7").Value = 2005
$ws.Range("B7").Value = "Brut Vintage"
$ws.Range("C7").Value = "Dom Pérignon"
$ws.Range("D7").Value = "Champagne"
$ws.Range("E7").Value = "Champagne"
$ws.Range("F7").Value = "France"
$ws.Range("G7").Value = 250
$ws.Range("H7").Value = 1100
$ws.Range("I7").Formula = "=G7/(H7/1.081)"
$ws.Range("J7").Value = 8
$ws.Range("K6").Copy($ws.Range("K7"))
$ws.Range("K7").Value = 45662
$ws.Range("L7").Value = 6

# ---- Row 8 : 2021 Patrick Piuze "Terroir Chapelle" (Chablis / Bourgogne) ----
$ws.Range("A8").Value = 2021
$ws.Range("B8").Value = "Terroir Chapelle"
$ws.Range("C8").Value = "Patrick Piuze"
$ws.Range("D8").Value = "Chablis"
$ws.Range("E8").Value = "Bourgogne"
$ws.Range("F8").Value = "France"
$ws.Range("G8").Value = 38
$ws.Range("H8").Value = 130
$ws.Range("I8").Formula = "=G8/(H8/1.081)"
$ws.Range("J8").Value = 12
$ws.Range("K6").Copy($ws.Range("K8"))
$ws.Range("K8").Value = 45662
$ws.Range("L8").Value = 2
$ws.Range("M8").Value = "Cavesa"

# Matches the author's final selection position after the edit.
[void]$ws.Range("C6").Select()
